$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Admin sheet: insert a new "Admin login" row above the existing
#    "Lấy danh sách quản trị viên" row (old row 4 -> new row 5).
# ---------------------------------------------------------------
$admin = $wb.Worksheets.Item("Admin")

# Insert a new row at position 4, pushing the old rows 4-6 down to 5-7.
$admin.Rows.Item(4).Insert()

$admin.Cells.Item(4, 1).Value = "Đăng nhập tài khoản quản trị viên"
$admin.Cells.Item(4, 2).Value = "/admins/login"
$admin.Cells.Item(4, 3).Value = "POST"
$admin.Cells.Item(4, 4).Value = "{
    username: String,
    password: String
}"
$admin.Cells.Item(4, 5).Value = "{
    error: false,
    message: ""Đăng nhập thành công"",
    data: {
        _id: String,
        username: String,
        password: String
    }
}"
$admin.Cells.Item(4, 6).Value = "{
    error: true,
    message: ""Tài khoản này không tồn tại""
}
hoặc
{
    error: true,
    message: ""Mật khẩu không chính xác""
}"

# match styles used by the rest of the data rows on this sheet
$admin.Cells.Item(4, 1).Style = $admin.Cells.Item(5, 1).Style
$admin.Cells.Item(4, 2).Style = $admin.Cells.Item(5, 2).Style
$admin.Cells.Item(4, 3).Style = $admin.Cells.Item(5, 3).Style
$admin.Cells.Item(4, 4).Style = $admin.Cells.Item(5, 4).Style
$admin.Cells.Item(4, 5).Style = $admin.Cells.Item(5, 5).Style
$admin.Cells.Item(4, 6).Style = $admin.Cells.Item(5, 6).Style

$admin.Rows.Item(4).RowHeight = 148.5

# ---------------------------------------------------------------
# 2. Update sheet view / selection on the "Admin" sheet.
# ---------------------------------------------------------------
$admin.Activate()
$admin.Range("A3").Select()
$admin.Application.ActiveWindow.ScrollRow = 3
$admin.Range("D3").Select()

# ---------------------------------------------------------------
# 3. Update sheet view / selection on the "User" sheet.
# ---------------------------------------------------------------
$user = $wb.Worksheets.Item("User")
$user.Activate()
$user.Range("A4").Select()
$user.Application.ActiveWindow.ScrollRow = 4
$user.Range("F4").Select()
